$wb = $excel.ActiveWorkbook

# ---- Sheet 2: "1_ISSUES" -> "SCHEME_MEASURES" ----
$s2 = $wb.Worksheets.Item("1_ISSUES")
# Remove the now-unused D:H columns (old "rule" sheet had 8 cols; new content only needs 3)
$s2.Range("D1:H2").Clear()
$s2.Cells.Item(1,1).Value = "Indicator"
$s2.Cells.Item(1,2).Value = "Description"
$s2.Cells.Item(1,3).Value = "Value"
$s2.Cells.Item(2,1).Value = "MQMS01"
$s2.Cells.Item(2,2).Value = "Total number of tables"
$s2.Cells.Item(2,3).Value = 5
$s2.Cells.Item(3,1).Value = "MQMS02"
$s2.Cells.Item(3,2).Value = "Total number of columns"
$s2.Cells.Item(3,3).Value = 46
$s2.Cells.Item(4,1).Value = "MQMS03"
$s2.Cells.Item(4,2).Value = "Total number of primary key"
$s2.Cells.Item(4,3).Value = 4
$s2.Cells.Item(5,1).Value = "MQMS04"
$s2.Cells.Item(5,2).Value = "Total number of foreign key"
$s2.Cells.Item(5,3).Value = 5
$s2.Cells.Item(6,1).Value = "MQMS05"
$s2.Cells.Item(6,2).Value = "Total number of unique key"
$s2.Cells.Item(6,3).Value = 0
$s2.Name = "SCHEME_MEASURES"

# ---- Sheet 3: "2_SCHEME_MEASURES" -> "METADATA_ISSUES" ----
$s3 = $wb.Worksheets.Item("2_SCHEME_MEASURES")
# Extend the bold/bordered header style from A1:C1 into the new D1:H1 header cells
$s3.Range("A1:C1").Copy()
$s3.Range("D1:H1").PasteSpecial(-4122)
# Remove the now-unused rows 3:6 (old scheme-measures sheet had 6 rows; new content only needs 2)
$s3.Range("A3:C6").Clear()
$s3.Cells.Item(1,1).Value = "rule"
$s3.Cells.Item(1,2).Value = "desc"
$s3.Cells.Item(1,3).Value = "owner"
$s3.Cells.Item(1,4).Value = "table"
$s3.Cells.Item(1,5).Value = "column"
$s3.Cells.Item(1,6).Value = "constraint_name"
$s3.Cells.Item(1,7).Value = "length"
$s3.Cells.Item(1,8).Value = "limit"
$s3.Cells.Item(2,1).Value = "MQME12"
$s3.Cells.Item(2,2).Value = "Total number of tables with non-standard column prefixes"
$s3.Cells.Item(2,3).Value = "SISAGUA"
$s3.Cells.Item(2,4).Value = "PEDIDO"
$s3.Cells.Item(2,5).Value = "UID_GERACAO_DEBITO"
$s3.Name = "METADATA_ISSUES"

# ---- Sheet 4: "3_MODEL_MEASURES" -> "METADATA_MEASURES" ----
$s4 = $wb.Worksheets.Item("3_MODEL_MEASURES")
$s4.Cells.Item(2,1).Value = "MQME00"
$s4.Cells.Item(2,2).Value = "Total number of columns"
$s4.Cells.Item(2,3).Value = 46
$s4.Cells.Item(3,1).Value = "MQMEA1"
$s4.Cells.Item(3,2).Value = "Total number of length-required columns"
$s4.Cells.Item(3,3).Value = 11
$s4.Cells.Item(4,1).Value = "MQMEA2"
$s4.Cells.Item(4,2).Value = "Total number of NUMBER columns"
$s4.Cells.Item(4,3).Value = 27
$s4.Name = "METADATA_MEASURES"

# ---- Sheet 5: "4_MODEL_METRICS" -> "METADATA_METRICS" ----
$s5 = $wb.Worksheets.Item("4_MODEL_METRICS")
# Remove the now-unused row 9 (old model-metrics sheet had 9 rows; new content only needs 8)
$s5.Range("A9:C9").Clear()
$s5.Cells.Item(2,1).Value = "IQME01"
$s5.Cells.Item(2,2).Value = "Columns with comments"
$s5.Cells.Item(2,3).NumberFormat = "@"
$s5.Cells.Item(2,3).Value = "100.00%"
$s5.Cells.Item(3,1).Value = "IQME02"
$s5.Cells.Item(3,2).Value = "Columns with data type"
$s5.Cells.Item(3,3).NumberFormat = "@"
$s5.Cells.Item(3,3).Value = "100.00%"
$s5.Cells.Item(4,1).Value = "IQME03"
$s5.Cells.Item(4,2).Value = "Length-required columns with data length"
$s5.Cells.Item(4,3).NumberFormat = "@"
$s5.Cells.Item(4,3).Value = "100.00%"
$s5.Cells.Item(5,1).Value = "IQME04"
$s5.Cells.Item(5,2).Value = "NUMBER columns with valid scale"
$s5.Cells.Item(5,3).NumberFormat = "@"
$s5.Cells.Item(5,3).Value = "100.00%"
$s5.Cells.Item(6,1).Value = "IQME05"
$s5.Cells.Item(6,2).Value = "Columns with valid num_distinct"
$s5.Cells.Item(6,3).NumberFormat = "@"
$s5.Cells.Item(6,3).Value = "100.00%"
$s5.Cells.Item(7,1).Value = "IQME06"
$s5.Cells.Item(7,2).Value = "Columns with valid num_nulls"
$s5.Cells.Item(7,3).NumberFormat = "@"
$s5.Cells.Item(7,3).Value = "100.00%"
$s5.Cells.Item(8,1).Value = "IQME07"
$s5.Cells.Item(8,2).Value = "Columns with valid density"
$s5.Cells.Item(8,3).NumberFormat = "@"
$s5.Cells.Item(8,3).Value = "100.00%"
$s5.Name = "METADATA_METRICS"
